# The paragraph about supported formats originally had the bold text split
# across two runs ("DOCX, DOC, PDF, HTML, XPS, R" + "TF and TXT") with a
# leftover "_GoBack" bookmark in between (an artifact Word leaves behind at
# the last edit position). Re-saving the document merges that text back into
# a single contiguous run and drops the stray bookmark.
$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Execute(
    "DOCX, DOC, PDF, HTML, XPS, RTF and TXT",  # FindText
    $false,                                    # MatchCase
    $false,                                    # MatchWholeWord
    $false,                                    # MatchWildcards
    $false,                                    # MatchSoundsLike
    $false,                                    # MatchAllWordForms
    $true,                                     # Forward
    1,                                         # Wrap (wdFindContinue)
    $false,                                    # Format
    "DOCX, DOC, PDF, HTML, XPS, RTF and TXT",  # ReplaceWith
    2                                          # Replace (wdReplaceAll)
)

$d.Save()
